$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 68.392882
$ws.Range("H2").Value = 205.178646
$ws.Range("I2").Value = 0.3817002623156464
$ws.Range("J2").Value = 0.3817002623156463
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 18.444833
$ws.Range("N2").Value = 55.33449900000001
$ws.Range("O2").Value = 0.529296397589589
$ws.Range("P2").Value = 0.5292963975895891
$ws.Range("Q2").Value = 1261.495286878706
$ws.Range("R2").Value = 11353.45758190835
$ws.Range("S2").Value = 0.2020325738026728
$ws.Range("T2").Value = 0.2020325738026728
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 68.392882
$ws.Range("H3").Value = 205.178646
$ws.Range("I3").Value = 0.3817002623156464
$ws.Range("J3").Value = 0.3817002623156463
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 8.028767999999999
$ws.Range("N3").Value = 24.086304
$ws.Range("O3").Value = 0.2303950368909585
$ws.Range("P3").Value = 0.2303950368909585
$ws.Range("Q3").Value = 549.110582429376
$ws.Range("R3").Value = 4941.995241864383
$ws.Range("S3").Value = 0.08794184601750188
$ws.Range("T3").Value = 0.08794184601750188
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 68.392882
$ws.Range("H4").Value = 205.178646
$ws.Range("I4").Value = 0.3817002623156464
$ws.Range("J4").Value = 0.3817002623156463
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 8.374233
$ws.Range("N4").Value = 25.122699
$ws.Range("O4").Value = 0.2403085655194523
$ws.Range("P4").Value = 0.2403085655194524
$ws.Range("Q4").Value = 572.737929409506
$ws.Range("R4").Value = 5154.641364685554
$ws.Range("S4").Value = 0.09172584249547165
$ws.Range("T4").Value = 0.09172584249547165
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 60.20577233333334
$ws.Range("H5").Value = 180.617317
$ws.Range("I5").Value = 0.3360080526004068
$ws.Range("J5").Value = 0.3360080526004068
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 18.444833
$ws.Range("N5").Value = 55.33449900000001
$ws.Range("O5").Value = 0.529296397589589
$ws.Range("P5").Value = 0.5292963975895891
$ws.Range("Q5").Value = 1110.485416324354
$ws.Range("R5").Value = 9994.368746919185
$ws.Range("S5").Value = 0.1778478518024885
$ws.Range("T5").Value = 0.1778478518024885
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 60.20577233333334
$ws.Range("H6").Value = 180.617317
$ws.Range("I6").Value = 0.3360080526004068
$ws.Range("J6").Value = 0.3360080526004068
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 8.028767999999999
$ws.Range("N6").Value = 24.086304
$ws.Range("O6").Value = 0.2303950368909585
$ws.Range("P6").Value = 0.2303950368909585
$ws.Range("Q6").Value = 483.378178325152
$ws.Range("R6").Value = 4350.403604926368
$ws.Range("S6").Value = 0.07741458767452986
$ws.Range("T6").Value = 0.07741458767452987
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 60.20577233333334
$ws.Range("H7").Value = 180.617317
$ws.Range("I7").Value = 0.3360080526004068
$ws.Range("J7").Value = 0.3360080526004068
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 8.374233
$ws.Range("N7").Value = 25.122699
$ws.Range("O7").Value = 0.2403085655194523
$ws.Range("P7").Value = 0.2403085655194524
$ws.Range("Q7").Value = 504.177165464287
$ws.Range("R7").Value = 4537.594489178584
$ws.Range("S7").Value = 0.08074561312338846
$ws.Range("T7").Value = 0.08074561312338846
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 50.58089766666667
$ws.Range("H8").Value = 151.742693
$ws.Range("I8").Value = 0.2822916850839468
$ws.Range("J8").Value = 0.2822916850839468
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 18.444833
$ws.Range("N8").Value = 55.33449900000001
$ws.Range("O8").Value = 0.529296397589589
$ws.Range("P8").Value = 0.5292963975895891
$ws.Range("Q8").Value = 932.9562104517564
$ws.Range("R8").Value = 8396.605894065808
$ws.Range("S8").Value = 0.1494159719844277
$ws.Range("T8").Value = 0.1494159719844278
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 50.58089766666667
$ws.Range("H9").Value = 151.742693
$ws.Range("I9").Value = 0.2822916850839468
$ws.Range("J9").Value = 0.2822916850839468
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 8.028767999999999
$ws.Range("N9").Value = 24.086304
$ws.Range("O9").Value = 0.2303950368909585
$ws.Range("P9").Value = 0.2303950368909585
$ws.Range("Q9").Value = 406.1022925974079
$ws.Range("R9").Value = 3654.920633376672
$ws.Range("S9").Value = 0.06503860319892675
$ws.Range("T9").Value = 0.06503860319892676
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 50.58089766666667
$ws.Range("H10").Value = 151.742693
$ws.Range("I10").Value = 0.2822916850839468
$ws.Range("J10").Value = 0.2822916850839468
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 8.374233
$ws.Range("N10").Value = 25.122699
$ws.Range("O10").Value = 0.2403085655194523
$ws.Range("P10").Value = 0.2403085655194524
$ws.Range("Q10").Value = 423.576222409823
$ws.Range("R10").Value = 3812.186001688407
$ws.Range("S10").Value = 0.06783710990059223
$ws.Range("T10").Value = 0.06783710990059223
